$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ColumnWidth = 26.85546875
$ws.Range("F1").Value = "ja"
